# Update PSP Time Recording Log sheet:
#  - F33 "node.js 5강" -> "node.js 8강" (effort description text correction)
#  - C33 stop time      0:42:59 AM -> ... (17:04 -> 17:30)
#  - E33 delta minutes   43 -> 68
#  - Row 34 (previously blank) filled in with a new log entry:
#       A34 = "11월 25일", B34 = 02:20 (start), C34 = 03:33 (stop),
#       D34 = 0 (interruption), E34 = 73 (delta minutes),
#       F34 = "요람 엑셀작성"
#  - Active selection moves from E34 to F35

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 edits -----------------------------------------------------
$ws.Range("C33").Value = 0.72916666666666663
$ws.Range("E33").Value = 68
$ws.Range("F33").Value = "node.js 8강"

# --- Row 34 edits (previously empty) -----------------------------------
$ws.Range("A34").Value = "11월 25일"
$ws.Range("B34").Value = 0.097222222222222224
$ws.Range("C34").Value = 0.14791666666666667
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 73
$ws.Range("F34").Value = "요람 엑셀작성"

# --- Selection update ---------------------------------------------------
$ws.Range("F35").Select()
